$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 82
$ws.Range("B82").Value = 5579144
$ws.Range("F82").Value = "Sabah"
$ws.Range("G82").Value = "Zira IK"
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = "D"
$ws.Range("K82").Value = 1.45
$ws.Range("L82").Value = 4.2
$ws.Range("M82").Value = 5.5
$ws.Range("N82").Value = 1.5
$ws.Range("O82").Value = 4
$ws.Range("P82").Value = 5.25
$ws.Range("Q82").Value = -1
$ws.Range("R82").Value = 1.85
$ws.Range("S82").Value = 1.95
$ws.Range("T82").Value = 2.5
$ws.Range("U82").Value = 1.8
$ws.Range("V82").Value = 2
$ws.Range("W82").Value = -1
$ws.Range("X82").Value = 3
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = -1
$ws.Range("AA82").Value = 0.95
$ws.Range("AB82").Value = -1
$ws.Range("AC82").Value = 1

# Row 84
$ws.Range("B84").Value = 5574442
$ws.Range("F84").Value = "FK Qarabag"
$ws.Range("G84").Value = "FK Sumqayit"
$ws.Range("H84").Value = 1
$ws.Range("I84").Value = 2
$ws.Range("J84").Value = "A"
$ws.Range("K84").Value = 1.125
$ws.Range("L84").Value = 7.5
$ws.Range("M84").Value = 15
$ws.Range("N84").Value = 1.2
$ws.Range("O84").Value = 6
$ws.Range("P84").Value = 11
$ws.Range("Q84").Value = -2.25
$ws.Range("R84").Value = 1.975
$ws.Range("S84").Value = 1.825
$ws.Range("T84").Value = 3.5
$ws.Range("U84").Value = 1.825
$ws.Range("V84").Value = 1.975
$ws.Range("W84").Value = -1
$ws.Range("X84").Value = -1
$ws.Range("Y84").Value = 10
$ws.Range("Z84").Value = -1
$ws.Range("AA84").Value = 0.825
$ws.Range("AB84").Value = -1
$ws.Range("AC84").Value = 0.9750000000000001

# Row 85
$ws.Range("B85").Value = 5573342
$ws.Range("F85").Value = "PFK Turan Tovuz"
$ws.Range("G85").Value = "Sabail FC"
$ws.Range("H85").Value = 2
$ws.Range("I85").Value = 2
$ws.Range("J85").Value = "D"
$ws.Range("K85").Value = 2.6
$ws.Range("L85").Value = 3
$ws.Range("M85").Value = 2.6
$ws.Range("N85").Value = 2.8
$ws.Range("O85").Value = 2.875
$ws.Range("P85").Value = 2.5
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = 2.05
$ws.Range("S85").Value = 1.75
$ws.Range("T85").Value = 2.25
$ws.Range("U85").Value = 1.875
$ws.Range("V85").Value = 1.925
$ws.Range("W85").Value = -1
$ws.Range("X85").Value = 1.875
$ws.Range("Y85").Value = -1
$ws.Range("Z85").Value = 0
$ws.Range("AA85").Value = -0
$ws.Range("AB85").Value = 0.875
$ws.Range("AC85").Value = -1

# Row 86
$ws.Range("B86").Value = 5602065
$ws.Range("F86").Value = "FK Gabala"
$ws.Range("G86").Value = "FK Kapaz"
$ws.Range("H86").Value = 1
$ws.Range("I86").Value = 1
$ws.Range("J86").Value = "D"
$ws.Range("K86").Value = 1.8
$ws.Range("L86").Value = 3.5
$ws.Range("M86").Value = 3.6
$ws.Range("N86").Value = 1.4
$ws.Range("O86").Value = 4.2
$ws.Range("P86").Value = 6
$ws.Range("Q86").Value = -1.25
$ws.Range("R86").Value = 1.9
$ws.Range("S86").Value = 1.9
$ws.Range("T86").Value = 2.5
$ws.Range("U86").Value = 1.85
$ws.Range("V86").Value = 1.95
$ws.Range("W86").Value = -1
$ws.Range("X86").Value = 3.2
$ws.Range("Y86").Value = -1
$ws.Range("Z86").Value = -1
$ws.Range("AA86").Value = 0.8999999999999999
$ws.Range("AB86").Value = -1
$ws.Range("AC86").Value = 0.95

# Row 87
$ws.Range("B87").Value = 5607916
$ws.Range("F87").Value = "FK Sumqayit"
$ws.Range("G87").Value = "Zira IK"
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 3
$ws.Range("J87").Value = "A"
$ws.Range("K87").Value = 4
$ws.Range("L87").Value = 3.8
$ws.Range("M87").Value = 1.666
$ws.Range("N87").Value = 4.333
$ws.Range("O87").Value = 3.6
$ws.Range("P87").Value = 1.65
$ws.Range("Q87").Value = 0.75
$ws.Range("R87").Value = 1.95
$ws.Range("S87").Value = 1.85
$ws.Range("T87").Value = 2.25
$ws.Range("U87").Value = 1.925
$ws.Range("V87").Value = 1.875
$ws.Range("W87").Value = -1
$ws.Range("X87").Value = -1
$ws.Range("Y87").Value = 0.6499999999999999
$ws.Range("Z87").Value = -1
$ws.Range("AA87").Value = 0.8500000000000001
$ws.Range("AB87").Value = 0.925
$ws.Range("AC87").Value = -1

